$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 812
$ws.Range("F3").Value = 555
$ws.Range("F5").Value = 509
$ws.Range("F6").Value = 1147
$ws.Range("F7").Value = 328
$ws.Range("F9").Value = 124
$ws.Range("F11").Value = 1189
$ws.Range("F12").Value = 53
$ws.Range("F14").Value = 861
$ws.Range("F17").Value = 64
$ws.Range("F20").Value = 754
$ws.Range("F22").Value = 2757
$ws.Range("F23").Value = 787
$ws.Range("F24").Value = 83
$ws.Range("F25").Value = 2101
$ws.Range("F26").Value = 663
$ws.Range("F27").Value = 2965
$ws.Range("F28").Value = 562
$ws.Range("F29").Value = 3
$ws.Range("F30").Value = 5
$ws.Range("F31").Value = 87
$ws.Range("F32").Value = 719
$ws.Range("F34").Value = 125
$ws.Range("F36").Value = 1042
$ws.Range("F37").Value = 1754
$ws.Range("F38").Value = 377
$ws.Range("F40").Value = 548
$ws.Range("F41").Value = 179
$ws.Range("F43").Value = 168
$ws.Range("F44").Value = 38

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 12
$ws.Range("F10").Value = 6
$ws.Range("F15").Value = 1

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 812
$ws.Range("F4").Value = 555
$ws.Range("F6").Value = 509
$ws.Range("F7").Value = 1147
$ws.Range("F8").Value = 328
$ws.Range("F10").Value = 124
$ws.Range("F12").Value = 1189
$ws.Range("F13").Value = 53
$ws.Range("F14").Value = 861
$ws.Range("F18").Value = 64
$ws.Range("F21").Value = 754
$ws.Range("F23").Value = 2757
$ws.Range("F24").Value = 787
$ws.Range("F25").Value = 83
$ws.Range("F28").Value = 2965
$ws.Range("F29").Value = 562
$ws.Range("F31").Value = 12
$ws.Range("F32").Value = 6
$ws.Range("F34").Value = 87
$ws.Range("F36").Value = 719
$ws.Range("F38").Value = 125
$ws.Range("F40").Value = 1042
$ws.Range("F41").Value = 1754
$ws.Range("F43").Value = 377
$ws.Range("F44").Value = 548
$ws.Range("F45").Value = 179
$ws.Range("F47").Value = 168
$ws.Range("F48").Value = 38

